$d = $word.ActiveDocument

# Step 1: Remove the empty paragraph right after the "aims to:" sentence.
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("aims to:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $rng1.Collapse(0)
    $emptyParaRange = $d.Range($rng1.End, $rng1.End + 1)
    $emptyParaRange.Delete()
}

# Step 2: Locate the "In doing so..." paragraph and replace/expand it with the
# new Background Material section. InsertXML at the collapsed end of this
# paragraph's text consumes the paragraph's trailing mark, so we re-emit the
# "In doing so" paragraph (without the lastRenderedPageBreak, which moves to
# the new heading) as the first paragraph of the inserted block. The leftover
# empty paragraph this produces at the very end of the body becomes the final
# blank paragraph the target document expects before the sectPr.
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("In doing so, this report seeks to provide insights and provoke thought on the effectiveness of current practices and the potential need for a paradigm shift in managing software development projects.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rng2.Collapse(0)
    $xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>In doing so, this report seeks to provide insights and provoke thought on the effectiveness of current practices and the potential need for a paradigm shift in managing software development projects.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:lastRenderedPageBreak/><w:t>BACKGROUND MATERIAL</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>To understand the importance of measuring the value of software development results, it is important to have a clear understanding of the context in which software development takes place. The background material necessary to understand the discussion in "Value Results, Not Just Effort" by Venkat Subramaniam encompasses several key subjects, including:</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>SUBJECT 1: The Importance of Extensibility in Software Development</w:t></w:r></w:p><w:p><w:r><w:t>Extensibility is a fundamental concept in software engineering that refers to the ability of a system to accommodate future growth and changes. It is about designing software in a way that new functionality can be added with minimal impact on the existing system. The principle of extensibility anticipates future needs and technological advancements, ensuring that the software remains relevant and adaptable over time.</w:t></w:r></w:p><w:p><w:r><w:t>However, as Subramaniam notes, the pursuit of extensibility can lead to the inclusion of superfluous code that may never be used, which can burden a project and distract from its immediate goals. Balancing the need for extensibility with the need for efficiency is crucial.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>SUBJECT 2: The Dangers of Overworking</w:t></w:r></w:p><w:p><w:r><w:t>Overworking in the software industry is a critical issue, often stemming from a culture that values long hours as a sign of commitment and productivity. This culture can lead to burnout, decreased job satisfaction, and diminished quality of work. The adverse effects of overworking are well-documented in occupational psychology. Prolonged periods of intense work without adequate rest can impair cognitive functions, reduce creativity, and negatively impact overall mental and physical health.</w:t></w:r></w:p><w:p><w:r><w:t>The concept of "work smarter, not harder" is increasingly being adopted as organizations recognize the limitations and drawbacks of overworking. This approach advocates for efficiency and effectiveness in work practices, emphasizing the importance of rest and work-life balance.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">SUBJECT 3: </w:t></w:r><w:r><w:t xml:space="preserve"> Evolution of Software Development Practices</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Software development has undergone significant transformations over the years. Initially dominated by waterfall methodologies, which are linear and sequential, the industry has increasingly moved towards agile methodologies. Agile practices, characterized by iterative development, flexibility, and collaboration, have revolutionized the way software is built and delivered. Understanding this evolution is critical to appreciating why measuring the value of results, as opposed to effort, has become increasingly relevant.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">SUBJECT </w:t></w:r><w:r><w:t>4</w:t></w:r><w:r><w:t xml:space="preserve">:  </w:t></w:r><w:r><w:t>The Challenges and Opportunities in Software Development</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Today''s software development is marked by a need to deliver products that not only meet technical requirements but also provide real value to users. This involves navigating challenges like rapidly changing technology, high customer expectations, and the necessity for user-</w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>centered design. The focus has shifted from just delivering software to ensuring that the software solves real problems and enhances user experiences.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">SUBJECT </w:t></w:r><w:r><w:t>5</w:t></w:r><w:r><w:t xml:space="preserve">:  </w:t></w:r><w:r><w:t>Current State of Software Development Practices</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Agile methodologies, continuous delivery, continuous integration, and user feedback are at the forefront of modern software development practices. These approaches emphasize the importance of adaptability, frequent delivery of working software, and close collaboration between developers and business stakeholders. Understanding these practices is vital for appreciating the need to value results over effort.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">SUBJECT </w:t></w:r><w:r><w:t>6</w:t></w:r><w:r><w:t xml:space="preserve">:  </w:t></w:r><w:r><w:t>Measuring the Value of Software Development Results</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>There is a growing body of research focused on how to effectively measure the value of software development efforts. These studies explore various approaches and metrics, highlight challenges in quantifying value, and point to the potential for future innovations in this area. A comprehensive understanding of this research is essential for anyone looking to understand the importance of valuing results in software development.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng2.InsertXML($xml2)
}
